$d = $word.ActiveDocument

# The "_GoBack" bookmark currently sits right after "Vratsa forever!!!",
# just before that paragraph's mark. We want to push two new paragraphs
# in after that text, and have the bookmark end up after all of the new
# text (inside the last new paragraph, after its run).
#
# InsertBefore() on a Range placed at the bookmark's point inserts text
# ahead of the bookmark; when the inserted text itself contains paragraph
# marks, the bookmark re-anchors to just after the last inserted mark, and
# plain (markless) text inserted the same way ends up before the bookmark.
# Doing this in two steps (paragraph breaks first, then the trailing
# sentence) leaves the bookmark exactly where Word would: at the very end
# of the new content.

$bm = $d.Bookmarks("_GoBack")
$p1 = $d.Range($bm.Start, $bm.Start)
$p1.InsertBefore("`r…`r")

$bm = $d.Bookmarks("_GoBack")
$p2 = $d.Range($bm.Start, $bm.Start)
$p2.InsertBefore("Now I want to add some more!")
